# Update "Pais" worksheet with latest COVID case data and timestamp
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Last-updated timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 15 de Septiembre de 2020 a las 00:05"

# Straightforward numeric updates (no re-ranking)
$ws.Range("B4").Value = 6744652
$ws.Range("C4").Value = 33582
$ws.Range("D4").Value = 4019427
$ws.Range("E4").Value = 2526354
$ws.Range("G4").Value = 351
$ws.Range("H4").Value = 198871
$ws.Range("B6").Value = 4345610
$ws.Range("C6").Value = 15155
$ws.Range("D6").Value = 3613184
$ws.Range("E6").Value = 600420
$ws.Range("G6").Value = 343
$ws.Range("H6").Value = 132006
$ws.Range("B8").Value = 733860
$ws.Range("C8").Value = 4241
$ws.Range("D8").Value = 573364
$ws.Range("E8").Value = 129684
$ws.Range("G8").Value = 102
$ws.Range("H8").Value = 30812
$ws.Range("B11").Value = 650749
$ws.Range("C11").Value = 956
$ws.Range("D11").Value = 579289
$ws.Range("E11").Value = 55961
$ws.Range("G11").Value = 52
$ws.Range("H11").Value = 15499
$ws.Range("B63").Value = 45493
$ws.Range("C63").Value = 59
$ws.Range("D63").Value = 44471
$ws.Range("E63").Value = 733
$ws.Range("G63").Value = 3
$ws.Range("H63").Value = 289
$ws.Range("B84").Value = 18061
$ws.Range("C84").Value = 143
$ws.Range("D84").Value = 12930
$ws.Range("E84").Value = 4402
$ws.Range("G84").Value = 9
$ws.Range("H84").Value = 729
$ws.Range("B100").Value = 9243
$ws.Range("C100").Value = 70
$ws.Range("D100").Value = 7536
$ws.Range("E100").Value = 1674
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = 33
$ws.Range("B102").Value = 8654
$ws.Range("C102").Value = 11
$ws.Range("D102").Value = 7785
$ws.Range("E102").Value = 816
$ws.Range("B115").Value = 5000
$ws.Range("C115").Value = 4
$ws.Range("D115").Value = 4496
$ws.Range("E115").Value = 421
$ws.Range("B122").Value = 4726
$ws.Range("C122").Value = 42
$ws.Range("D122").Value = 4040
$ws.Range("E122").Value = 578
$ws.Range("B151").Value = 2111
$ws.Range("C151").Value = 2
$ws.Range("E151").Value = 403
$ws.Range("B152").Value = 2013
$ws.Range("C152").Value = 2
$ws.Range("D152").Value = 1215
$ws.Range("E152").Value = 215
$ws.Range("B153").Value = 1884
$ws.Range("C153").Value = 31
$ws.Range("D153").Value = 1265
$ws.Range("E153").Value = 563
$ws.Range("B156").Value = 1717
$ws.Range("C156").Value = 10
$ws.Range("D156").Value = 1137
$ws.Range("E156").Value = 524
$ws.Range("B165").Value = 1085
$ws.Range("C165").Value = 1
$ws.Range("D165").Value = 940
$ws.Range("E165").Value = 64
$ws.Range("G165").Value = 1
$ws.Range("H165").Value = 81
$ws.Range("B189").Value = 183
$ws.Range("C189").Value = 2
$ws.Range("D189").Value = 165
$ws.Range("E189").Value = 11

# Botsuana overtakes Malta/Georgia/Guinea-Bisau/Benin in total cases,
# shifting rows 145-149 down; row 145 becomes Botsuana with fresh data
# while Malta, Georgia, Guinea-Bisau, Benin shift to rows 146-149.
$ws.Range("A145").Value = "Botsuana"
$ws.Range("B145").Value = 2463
$ws.Range("C145").Value = 211
$ws.Range("D145").Value = 575
$ws.Range("E145").Value = 1877
$ws.Range("H145").Value = 11

$ws.Range("A146").Value = "Malta"
$ws.Range("B146").Value = 2405
$ws.Range("C146").Value = 53
$ws.Range("D146").Value = 1890
$ws.Range("E146").Value = 499
$ws.Range("G146").Value = 1
$ws.Range("H146").Value = 16

$ws.Range("A147").Value = "Georgia"
$ws.Range("B147").Value = 2392
$ws.Range("C147").Value = 165
$ws.Range("D147").Value = 1369
$ws.Range("E147").Value = 1004
$ws.Range("H147").Value = 19

$ws.Range("A148").Value = "Guinea-Bisau"
$ws.Range("B148").Value = 2275
$ws.Range("D148").Value = 1127
$ws.Range("E148").Value = 1109
$ws.Range("H148").Value = 39

$ws.Range("A149").Value = "Benin"
$ws.Range("B149").Value = 2267
$ws.Range("D149").Value = 1942
$ws.Range("E149").Value = 285
$ws.Range("H149").Value = 40

